$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Stamp the "About" sheet with the date the template was last refreshed,
# placed in C1 next to the title in A1, formatted as a date.
# Excel serial date 44307 == 2021-04-21.
$ws.Range("C1").Value = 44307
$ws.Range("C1").NumberFormat = "mm-dd-yy"
